{"js": "// Update the date paragraph and every arithmetic-problem cell in the single\n// table, preserving existing run formatting (font, size, etc.) by replacing\n// only the text content (same approach as the underlying OOXML diff, which\n// only touches each <w:t> element's text).\n\n// --- 1. Update the title/date paragraph -----------------------------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateOld = \"2024-06-28 Friday\";\nconst dateNew = \"2024-06-29 Saturday\";\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\n\nif (titlePara.text.trim() === dateOld) {\n  titlePara.getRange().insertText(dateNew, Word.InsertLocation.replace);\n}\n\n// --- 2. Update every cell value in the practice table ----------------------\n// Old -> new values, in the same row-major order as the table (20 rows x 5\n// columns = 100 cells), matching the order cells appear in the document.\nconst newValues = [\n  [\"31-20=11\", \"66-36=30\", \"91-42=49\", \"23+61=84\", \"83-34=49\"],\n  [\"57-37=20\", \"21-1=20\", \"87-11=76\", \"80+16=96\", \"65+21=86\"],\n  [\"77-19=58\", \"96-81=15\", \"11+88=99\", \"6+22=28\", \"80-15=65\"],\n  [\"49+23=72\", \"25+18=43\", \"74-48=26\", \"68-19=49\", \"66-65=1\"],\n  [\"82-31=51\", \"51-24=27\", \"34-19=15\", \"49+44=93\", \"39-38=1\"],\n  [\"92-82=10\", \"42+9=51\", \"26+70=96\", \"58-21=37\", \"85-35=50\"],\n  [\"44-16=28\", \"60+13=73\", \"58-30=28\", \"50+30=80\", \"64-40=24\"],\n  [\"36+53=89\", \"19-8=11\", \"20+59=79\", \"48-34=14\", \"67-0=67\"],\n  [\"8+64=72\", \"59-14=45\", \"58-44=14\", \"17+12=29\", \"77-56=21\"],\n  [\"33+31=64\", \"28-2=26\", \"81-43=38\", \"29+21=50\", \"72-45=27\"],\n  [\"47+16=63\", \"23-20=3\", \"13-8=5\", \"31-26=5\", \"87-28=59\"],\n  [\"69-11=58\", \"41-17=24\", \"66+28=94\", \"49+30=79\", \"36+5=41\"],\n  [\"26+64=90\", \"85+4=89\", \"99-86=13\", \"70-64=6\", \"69+18=87\"],\n  [\"96-95=1\", \"1+8=9\", \"91-60=31\", \"6-4=2\", \"89-56=33\"],\n  [\"46+36=82\", \"44+44=88\", \"60+5=65\", \"15+12=27\", \"76-22=54\"],\n  [\"76-70=6\", \"54-10=44\", \"11+44=55\", \"65-42=23\", \"5+8=13\"],\n  [\"24+43=67\", \"91-24=67\", \"31+0=31\", \"54-20=34\", \"48-42=6\"],\n  [\"60-11=49\", \"8+77=85\", \"6+15=21\", \"36-2=34\", \"79-35=44\"],\n  [\"29+65=94\", \"53-16=37\", \"2+22=24\", \"25+26=51\", \"49+9=58\"],\n  [\"89-87=2\", \"67+21=88\", \"28+7=35\", \"66+8=74\", \"11+34=45\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the date paragraph and every arithmetic-problem cell in the single\n# table, preserving existing run formatting by assigning Range.Text (which\n# replaces the text content in place, leaving the run's rPr/pPr untouched),\n# mirroring the OOXML diff which only rewrites each <w:t> element's text.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the title/date paragraph ------------------------------------\n$dateOld = \"2024-06-28 Friday\"\n$dateNew = \"2024-06-29 Saturday\"\n\n$titlePara = $d.Paragraphs(1)\nif ($titlePara.Range.Text.TrimEnd(\"`r\", \"`n\") -eq $dateOld) {\n    $titlePara.Range.Text = $dateNew\n}\n\n# --- 2. Update every cell value in the practice table -----------------------\n# New values, row-major (20 rows x 5 columns), matching Table.Cell(row, col)\n# iteration order and the order cells appear in the document.\n$newValues = @(\n  @(\"31-20=11\", \"66-36=30\", \"91-42=49\", \"23+61=84\", \"83-34=49\"),\n  @(\"57-37=20\", \"21-1=20\", \"87-11=76\", \"80+16=96\", \"65+21=86\"),\n  @(\"77-19=58\", \"96-81=15\", \"11+88=99\", \"6+22=28\", \"80-15=65\"),\n  @(\"49+23=72\", \"25+18=43\", \"74-48=26\", \"68-19=49\", \"66-65=1\"),\n  @(\"82-31=51\", \"51-24=27\", \"34-19=15\", \"49+44=93\", \"39-38=1\"),\n  @(\"92-82=10\", \"42+9=51\", \"26+70=96\", \"58-21=37\", \"85-35=50\"),\n  @(\"44-16=28\", \"60+13=73\", \"58-30=28\", \"50+30=80\", \"64-40=24\"),\n  @(\"36+53=89\", \"19-8=11\", \"20+59=79\", \"48-34=14\", \"67-0=67\"),\n  @(\"8+64=72\", \"59-14=45\", \"58-44=14\", \"17+12=29\", \"77-56=21\"),\n  @(\"33+31=64\", \"28-2=26\", \"81-43=38\", \"29+21=50\", \"72-45=27\"),\n  @(\"47+16=63\", \"23-20=3\", \"13-8=5\", \"31-26=5\", \"87-28=59\"),\n  @(\"69-11=58\", \"41-17=24\", \"66+28=94\", \"49+30=79\", \"36+5=41\"),\n  @(\"26+64=90\", \"85+4=89\", \"99-86=13\", \"70-64=6\", \"69+18=87\"),\n  @(\"96-95=1\", \"1+8=9\", \"91-60=31\", \"6-4=2\", \"89-56=33\"),\n  @(\"46+36=82\", \"44+44=88\", \"60+5=65\", \"15+12=27\", \"76-22=54\"),\n  @(\"76-70=6\", \"54-10=44\", \"11+44=55\", \"65-42=23\", \"5+8=13\"),\n  @(\"24+43=67\", \"91-24=67\", \"31+0=31\", \"54-20=34\", \"48-42=6\"),\n  @(\"60-11=49\", \"8+77=85\", \"6+15=21\", \"36-2=34\", \"79-35=44\"),\n  @(\"29+65=94\", \"53-16=37\", \"2+22=24\", \"25+26=51\", \"49+9=58\"),\n  @(\"89-87=2\", \"67+21=88\", \"28+7=35\", \"66+8=74\", \"11+34=45\")\n)\n\n$table = $d.Tables(1)\nfor ($r = 1; $r -le 20; $r++) {\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cell.Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
